$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule R30's lower bound (C10) is updated from 18 to 1
$ws.Range("C10").Value = 1
